$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOB1211")

# Update the course name text (row 3, columns B and C)
$ws.Range("B3").Value = " Poluição Atmosférica"
$ws.Range("C3").Value = " Poluição Atmosférica"

# Update the "Semestre ideal" value (row 9, columns B and C)
$ws.Range("B9").Value = "EA-7"
$ws.Range("C9").Value = "EA-7"

# Remove the "Requisitos" rows (23-26), which are no longer part of this sheet
$ws.Range("A23:C26").EntireRow.Delete()
